# Better Character Tooltips - Use bright (instead of ordinary) yellow for
# medium-level wealth & artifacts.
#
# Sheet "Artifact_Keys" has a small lookup table in E2:G10 that maps a
# wealth/artifact quality tier (1-9) to a colour-code prefix/suffix used to
# build the generated localisation text in column G (rows 12-56).
#
# F4 holds the colour code used for tier 3 (Q3): ordinary yellow "$YQ" ->
# change to medium/bright yellow "$MQ".
# F7/F8 hold the colour code used for tiers 6 and 7 (Q6/Q7): they currently
# share the same "medium" yellow "$MQ" that F4 used to have - bump them up
# to the brighter yellow "$lQ" (matching tiers 8/9 in F9/F10).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artifact_Keys")
$ws.Activate()

$ws.Range("F4").Value = [string][char]167 + "MQ"
$ws.Range("F7").Value = [string][char]167 + "lQ"
$ws.Range("F8").Value = [string][char]167 + "lQ"

# Match the recorded view-state change: selection moved from A1/F4 to B1/B7.
$ws.Range("B1").Select()
$ws.Range("B7").Select()
